# Auto-generated script applying scheduled market-data refresh to Adamantoise_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1246.3334
$ws.Range("I19").Value = 2499
$ws.Range("K19").Value = 2499
$ws.Range("M19").Value = -2324

$ws.Range("H74").Value = 6776.6
$ws.Range("I74").Value = 10053.2
$ws.Range("K74").Value = 10053.2
$ws.Range("M74").Value = -9117.200000000001

$ws.Range("H76").Value = 4786
$ws.Range("I76").Value = 4782.3335
$ws.Range("J76").Value = 4789.143
$ws.Range("K76").Value = 4782.3335
$ws.Range("L76").Value = 4789.143
$ws.Range("M76").Value = -4467.3335
$ws.Range("N76").Value = -5419.143

$ws.Range("H77").Value = 6776.6
$ws.Range("I77").Value = 10053.2
$ws.Range("K77").Value = 50266
$ws.Range("M77").Value = -45586

$ws.Range("H79").Value = 4786
$ws.Range("I79").Value = 4782.3335
$ws.Range("J79").Value = 4789.143
$ws.Range("K79").Value = 4782.3335
$ws.Range("L79").Value = 4789.143
$ws.Range("M79").Value = -3690.3335
$ws.Range("N79").Value = -6973.143

$ws.Range("H98").Value = 1390.55
$ws.Range("I98").Value = 1361.2354
$ws.Range("K98").Value = 1361.2354
$ws.Range("M98").Value = 136.7646

$ws.Range("H100").Value = 2715.4736
$ws.Range("I100").Value = 1544.125
$ws.Range("K100").Value = 1544.125
$ws.Range("M100").Value = -1003.125

$ws.Range("H105").Value = 52000
$ws.Range("J105").Value = 52000
$ws.Range("L105").Value = 52000
$ws.Range("N105").Value = -58988

$ws.Range("H106").Value = 4168791.2
$ws.Range("I106").Value = 4168791.2
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 4168791.2
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -4168160.2
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value = 3939
$ws.Range("I107").Value = 4479
$ws.Range("K107").Value = 4479
$ws.Range("M107").Value = -2559

$ws.Range("H122").Value = 1390.55
$ws.Range("I122").Value = 1361.2354
$ws.Range("K122").Value = 4083.7062
$ws.Range("M122").Value = -1633.7062

$ws.Range("H128").Value = 141900
$ws.Range("J128").Value = 141900
$ws.Range("L128").Value = 141900
$ws.Range("N128").Value = -151860

$ws.Range("H130").Value = 116995.5
$ws.Range("J130").Value = 116995.5
$ws.Range("L130").Value = 116995.5
$ws.Range("N130").Value = -127035.5

$ws.Range("H132").Value = 1706.7894
$ws.Range("I132").Value = 1706.7894
$ws.Range("K132").Value = 5120.3682
$ws.Range("M132").Value = -2590.3682

$ws.Range("H135").Value = 1562.5
$ws.Range("I135").Value = 1562.5
$ws.Range("K135").Value = 14062.5
$ws.Range("M135").Value = -11527.5

$ws.Range("H137").Value = 5559002
$ws.Range("I137").Value = 3804
$ws.Range("K137").Value = 11412
$ws.Range("M137").Value = -8862

$ws.Range("H138").Value = 2531.5557
$ws.Range("I138").Value = 1294
$ws.Range("J138").Value = 2864.7437
$ws.Range("K138").Value = 3882
$ws.Range("L138").Value = 8594.231100000001
$ws.Range("M138").Value = 1258
$ws.Range("N138").Value = -18874.2311

$ws.Range("H141").Value = 5199.4116
$ws.Range("I141").Value = 5199.4116
$ws.Range("K141").Value = 15598.2348
$ws.Range("M141").Value = -10418.2348


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2298.4614
$ws.Range("I2").Value = 890
$ws.Range("J2").Value = 2554.5454
$ws.Range("K2").Value = 890
$ws.Range("L2").Value = 2554.5454
$ws.Range("M2").Value = -777
$ws.Range("N2").Value = -2780.5454

$ws.Range("H32").Value = 24315916
$ws.Range("I32").Value = 25569254
$ws.Range("K32").Value = 25569254
$ws.Range("M32").Value = -25568967

$ws.Range("H116").Value = 2298.4614
$ws.Range("I116").Value = 890
$ws.Range("J116").Value = 2554.5454
$ws.Range("K116").Value = 890
$ws.Range("L116").Value = 2554.5454
$ws.Range("M116").Value = 1404
$ws.Range("N116").Value = -7142.5454

$ws.Range("H132").Value = 2584.1667
$ws.Range("I132").Value = 2226.6667
$ws.Range("K132").Value = 6680.000100000001
$ws.Range("M132").Value = -4150.000100000001

$ws.Range("H133").Value = 58000
$ws.Range("J133").Value = 58000
$ws.Range("L133").Value = 58000
$ws.Range("N133").Value = -63060

$ws.Range("H135").Value = 64184.8
$ws.Range("J135").Value = 64184.8
$ws.Range("L135").Value = 64184.8
$ws.Range("N135").Value = -74324.8


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2298.4614
$ws.Range("I3").Value = 890
$ws.Range("J3").Value = 2554.5454
$ws.Range("K3").Value = 890
$ws.Range("L3").Value = 2554.5454
$ws.Range("M3").Value = -776
$ws.Range("N3").Value = -2782.5454

$ws.Range("H20").Value = 35609.465
$ws.Range("I20").Value = 43311.918
$ws.Range("J20").Value = 4799.6665
$ws.Range("K20").Value = 43311.918
$ws.Range("L20").Value = 4799.6665
$ws.Range("M20").Value = -43064.918
$ws.Range("N20").Value = -5293.6665

$ws.Range("H105").Value = 2645.6155
$ws.Range("J105").Value = 2121.75
$ws.Range("L105").Value = 2121.75
$ws.Range("N105").Value = -5615.75


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8741.75
$ws.Range("I31").Value = 2796.75
$ws.Range("J31").Value = 14686.75
$ws.Range("K31").Value = 2796.75
$ws.Range("L31").Value = 14686.75
$ws.Range("M31").Value = -2501.75
$ws.Range("N31").Value = -15276.75

$ws.Range("H34").Value = 8741.75
$ws.Range("I34").Value = 2796.75
$ws.Range("J34").Value = 14686.75
$ws.Range("K34").Value = 2796.75
$ws.Range("L34").Value = 14686.75
$ws.Range("M34").Value = -2594.75
$ws.Range("N34").Value = -15090.75

$ws.Range("H62").Value = 4562
$ws.Range("I62").Value = 4224
$ws.Range("J62").Value = 4900
$ws.Range("K62").Value = 4224
$ws.Range("L62").Value = 4900
$ws.Range("M62").Value = -3600
$ws.Range("N62").Value = -6148

$ws.Range("H65").Value = 4562
$ws.Range("I65").Value = 4224
$ws.Range("J65").Value = 4900
$ws.Range("K65").Value = 21120
$ws.Range("L65").Value = 24500
$ws.Range("M65").Value = -18000
$ws.Range("N65").Value = -30740

$ws.Range("H132").Value = 1984.0555
$ws.Range("I132").Value = 1773.931
$ws.Range("K132").Value = 5321.793
$ws.Range("M132").Value = -2791.793


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2500250
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 100
$ws.Range("K26").Value = 300
$ws.Range("M26").Value = -12

$ws.Range("H62").Value = 3571.4285
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3571.4285
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 10714.2855
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -12086.2855

$ws.Range("H65").Value = 3571.4285
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3571.4285
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 32142.8565
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -39006.8565

$ws.Range("H124").Value = 5033
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 432.2353
$ws.Range("I2").Value = 390.33334
$ws.Range("K2").Value = 390.33334
$ws.Range("M2").Value = -277.33334

$ws.Range("H70").Value = 16624.75
$ws.Range("I70").Value = 33649.4
$ws.Range("K70").Value = 33649.4
$ws.Range("M70").Value = -33379.4

$ws.Range("H73").Value = 16624.75
$ws.Range("I73").Value = 33649.4
$ws.Range("K73").Value = 33649.4
$ws.Range("M73").Value = -32713.4

$ws.Range("H107").Value = 1250
$ws.Range("I107").Value = 1250
$ws.Range("K107").Value = 1250
$ws.Range("M107").Value = 670

$ws.Range("H113").Value = 17799.75
$ws.Range("J113").Value = 51492.25
$ws.Range("L113").Value = 51492.25
$ws.Range("N113").Value = -55832.25

$ws.Range("H119").Value = 41500.5
$ws.Range("J119").Value = 41500.5
$ws.Range("L119").Value = 41500.5
$ws.Range("N119").Value = -51176.5

$ws.Range("H122").Value = 1910.4667
$ws.Range("I122").Value = 1765.1666
$ws.Range("K122").Value = 5295.4998
$ws.Range("M122").Value = -2845.4998

$ws.Range("H132").Value = 2738.3572
$ws.Range("I132").Value = 2738.3572
$ws.Range("K132").Value = 8215.071599999999
$ws.Range("M132").Value = -5685.071599999999

$ws.Range("H136").Value = 62974.375
$ws.Range("J136").Value = 62974.375
$ws.Range("L136").Value = 188923.125
$ws.Range("N136").Value = -194023.125


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

$ws.Range("H133").Value = 29858
$ws.Range("J133").Value = 29858
$ws.Range("L133").Value = 29858
$ws.Range("N133").Value = -34918


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4974.2974
$ws.Range("I122").Value = 5661.8
$ws.Range("J122").Value = 3542
$ws.Range("K122").Value = 16985.4
$ws.Range("L122").Value = 10626
$ws.Range("M122").Value = -14535.4
$ws.Range("N122").Value = -15526

$ws.Range("H132").Value = 3117.6775
$ws.Range("I132").Value = 2842.4443
$ws.Range("J132").Value = 2842.4443
$ws.Range("K132").Value = 8527.332900000001
$ws.Range("M132").Value = -5997.332900000001

